# Reproduce the "Add files via upload" commit:
#   - add a new shared string "Copyright Infosys"
#   - write it into cell E31 of Sheet1
#   - update the view: scroll so row 26 is at the top, select E32

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New cell content (this also grows the shared-string table and the
# sheet's used-range/dimension automatically).
$ws.Range("E31").Value = "Copyright Infosys"

# Update the window/selection state to match the saved view.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 26
$win.ScrollColumn = 1
$ws.Range("E32").Select() | Out-Null
